# Avances hasta el 08-06-2019
#
# Schedule sheet "Hoja1": the Jueves (F10) and Viernes (F11) observation
# codes move from the "EE" series to the "EF" series.
#   F10: EE3 -> EF2
#   F11: EE4 -> EF1
# Write F11 first so the shared-string table gets the new unique strings
# appended in the same order the source workbook has them (02IE, EF1, EF2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F11").Value = "EF1"
$ws.Range("F10").Value = "EF2"

# The author's last action before saving left the active cell on F11.
[void]$ws.Range("F11").Select()
